# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Kujata_Profits workbook. For each affected Leve row, the price/profit
# columns (H:N) are refreshed with newly polled Market Board data.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 19
$ws.Range("H19").Value = 2569.3
$ws.Range("I19").Value = 3600
$ws.Range("J19").Value = 2127.5715
$ws.Range("K19").Value = 3600
$ws.Range("L19").Value = 2127.5715
$ws.Range("M19").Value = -3425
$ws.Range("N19").Value = -2477.5715

# Row 111
$ws.Range("H111").Value = 1204.2354
$ws.Range("I111").Value = 719.2222
$ws.Range("J111").Value = 1749.875
$ws.Range("K111").Value = 2157.6666
$ws.Range("L111").Value = 5249.625
$ws.Range("M111").Value = 909.3334
$ws.Range("N111").Value = -11383.625

# Row 138
$ws.Range("H138").Value = 751494.4399999999
$ws.Range("J138").Value = 1117036.4
$ws.Range("L138").Value = 3351109.2
$ws.Range("N138").Value = -3361389.2


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 45
$ws.Range("H45").Value = 1930.5714
$ws.Range("I45").Value = 1835.6666
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 1835.6666
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -1458.6666
$ws.Range("N45").Value = -3254

# Row 97
$ws.Range("H97").Value = 425.64285
$ws.Range("I97").Value = 429
$ws.Range("J97").Value = 405.5
$ws.Range("K97").Value = 429
$ws.Range("L97").Value = 405.5
$ws.Range("M97").Value = 67
$ws.Range("N97").Value = -1397.5


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 81
$ws.Range("H81").Value = 7053.909
$ws.Range("J81").Value = 7053.909
$ws.Range("L81").Value = 7053.909
$ws.Range("N81").Value = -9175.909

# Row 84
$ws.Range("H84").Value = 7053.909
$ws.Range("J84").Value = 7053.909
$ws.Range("L84").Value = 21161.727
$ws.Range("N84").Value = -31769.727

# Row 99
$ws.Range("H99").Value = 23810708
$ws.Range("I99").Value = 26316936
$ws.Range("K99").Value = 26316936
$ws.Range("M99").Value = -26315438

# Row 107
$ws.Range("H107").Value = 1471.0625
$ws.Range("I107").Value = 1156.7273
$ws.Range("J107").Value = 2162.6
$ws.Range("K107").Value = 1156.7273
$ws.Range("L107").Value = 2162.6
$ws.Range("M107").Value = 763.2727
$ws.Range("N107").Value = -6002.6


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 888.72546
$ws.Range("I31").Value = 746.0714
$ws.Range("J31").Value = 1554.4445
$ws.Range("K31").Value = 746.0714
$ws.Range("L31").Value = 1554.4445
$ws.Range("M31").Value = -451.0714
$ws.Range("N31").Value = -2144.4445

# Row 34
$ws.Range("H34").Value = 888.72546
$ws.Range("I34").Value = 746.0714
$ws.Range("J34").Value = 1554.4445
$ws.Range("K34").Value = 746.0714
$ws.Range("L34").Value = 1554.4445
$ws.Range("M34").Value = -544.0714
$ws.Range("N34").Value = -1958.4445


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 31
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("N31").ClearContents()

# Row 35
$ws.Range("H35").Value = 2750.75
$ws.Range("J35").Value = 4001.5
$ws.Range("L35").Value = 12004.5
$ws.Range("N35").Value = -12580.5

# Row 49
$ws.Range("H49").Value = 2001.3334
$ws.Range("J49").Value = 2001.3334
$ws.Range("L49").Value = 6004.0002
$ws.Range("N49").Value = -6316.0002

# Row 54
$ws.Range("H54").Value = 4504
$ws.Range("J54").Value = 4504
$ws.Range("L54").Value = 13512
$ws.Range("N54").Value = -14630

# Row 57
$ws.Range("H57").Value = 555
$ws.Range("I57").Value = 555
$ws.Range("K57").Value = 1665
$ws.Range("M57").Value = -1106

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("N77").ClearContents()

# Row 94
$ws.Range("H94").Value = 5354.5454
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 5354.5454
$ws.Range("K94").Value = 0
$ws.Range("N94").Value = -17415.6362
$ws.Range("M94").ClearContents()

# Row 99
$ws.Range("H99").Value = 2347.4285
$ws.Range("I99").Value = 674
$ws.Range("K99").Value = 2022
$ws.Range("M99").Value = 224

# Row 100
$ws.Range("H100").Value = 3426.8235
$ws.Range("J100").Value = 3426.8235
$ws.Range("L100").Value = 10280.4705
$ws.Range("N100").Value = -11902.4705

# Row 101
$ws.Range("H101").Value = 3989
$ws.Range("J101").Value = 3989
$ws.Range("L101").Value = 11967
$ws.Range("N101").Value = -16835

# Row 114
$ws.Range("H114").Value = 700.0526
$ws.Range("I114").Value = 395.6
$ws.Range("J114").Value = 1038.3334
$ws.Range("K114").Value = 1186.8
$ws.Range("L114").Value = 3115.0002
$ws.Range("M114").Value = 2067.2
$ws.Range("N114").Value = -9623.0002

# Row 129
$ws.Range("H129").Value = 19842410
$ws.Range("I129").Value = 37037884
$ws.Range("J129").Value = 6945803
$ws.Range("K129").Value = 111113652
$ws.Range("L129").Value = 20837409
$ws.Range("M129").Value = -111108652
$ws.Range("N129").Value = -20847409


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 1485.8518
$ws.Range("I102").Value = 1550.2222
$ws.Range("J102").Value = 1357.1111
$ws.Range("K102").Value = 1550.2222
$ws.Range("L102").Value = 1357.1111
$ws.Range("M102").Value = 71.77780000000007
$ws.Range("N102").Value = -4601.1111

# Row 104
$ws.Range("H104").Value = 68333.336
$ws.Range("J104").Value = 68333.336
$ws.Range("L104").Value = 68333.336
$ws.Range("N104").Value = -75321.336

# Row 107
$ws.Range("H107").Value = 606.0714
$ws.Range("I107").Value = 452.5
$ws.Range("J107").Value = 990
$ws.Range("K107").Value = 452.5
$ws.Range("L107").Value = 990
$ws.Range("M107").Value = 1467.5
$ws.Range("N107").Value = -4830

# Row 122
$ws.Range("H122").Value = 3025.5715
$ws.Range("I122").Value = 2382.375
$ws.Range("J122").Value = 3883.1667
$ws.Range("K122").Value = 7147.125
$ws.Range("L122").Value = 11649.5001
$ws.Range("M122").Value = -4697.125
$ws.Range("N122").Value = -16549.5001


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 112
$ws.Range("H112").Value = 47998.832
$ws.Range("J112").Value = 47998.832
$ws.Range("L112").Value = 47998.832
$ws.Range("N112").Value = -50952.832

# Row 122
$ws.Range("H122").Value = 31252114
$ws.Range("I122").Value = 31252114
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 93756342
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -93753892
$ws.Range("M122").ClearContents()


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N90").ClearContents()

# Row 122
$ws.Range("H122").Value = 20002370
$ws.Range("I122").Value = 21668984
$ws.Range("K122").Value = 65006952
$ws.Range("M122").Value = -65004502

# Row 132
$ws.Range("H132").Value = 2458.5908
$ws.Range("I132").Value = 2285.8667
$ws.Range("J132").Value = 2828.7144
$ws.Range("K132").Value = 6857.6001
$ws.Range("L132").Value = 8486.143199999999
$ws.Range("M132").Value = -4327.6001
$ws.Range("N132").Value = -13546.1432

# Row 136
$ws.Range("H136").Value = 2028.1428
$ws.Range("I136").Value = 1839.6
$ws.Range("J136").Value = 2499.5
$ws.Range("K136").Value = 5518.799999999999
$ws.Range("L136").Value = 7498.5
$ws.Range("M136").Value = -2968.799999999999
$ws.Range("N136").Value = -12598.5

